$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 18) holds a "Förändrad" date that was bumped
# by one day (45184 -> 45185, i.e. 2023-09-15 -> 2023-09-16).
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
